$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.828.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "'2.490.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'532.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'135.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").Value = "'5.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "'0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "'2.935.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "'58.738.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "'22.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'2.503.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "'10.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "'321.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'5.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("E23").Value = "  +4.18%  "
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "'7.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'0.0₃0754"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "'171.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  +4.59%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'18.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").Value = "'3.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").Value = "'3.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'280.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").Value = "'0.786"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D41").Value = "'5.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +2.94%  "
$ws.Range("D44").Value = "'129.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.97%  "
$ws.Range("D45").Value = "'10.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "'0.0920"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "'0.0496"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'17.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "'1.750.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "'0.981"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
